# Update the cryptocurrency price (column D) and 1h volume change (column E)
# figures on the active worksheet, matching the latest scrape performed by
# the "Updated cryptos list" GitHub Actions workflow.
#
# Column D values are formatted as text (e.g. "26.218.00", "0.0\u20858200")
# rather than genuine numbers, so each D cell's NumberFormat is forced to
# "@" (Text) before the new value is written -- this stops Excel's
# automatic type inference from reinterpreting strings like "139.00" or
# "1.280" as numbers and silently dropping the trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.218.00"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.660.03"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.48"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5219"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06333"
$ws.Range("E9").Value = "  -0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.03"
$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.653.62"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.425"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.889.18"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5466"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8200"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.96"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.251.87"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.641"
$ws.Range("E20").Value = "  -2.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.17"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.13"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.076"
$ws.Range("E23").Value = "  -4.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.00"
$ws.Range("E25").Value = "  -3.18%  "

$ws.Range("E26").Value = "  -3.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.229"
$ws.Range("E27").Value = "  -2.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.17"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05948"
$ws.Range("E30").Value = "  -3.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.280"
$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.624"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.300"
$ws.Range("E33").Value = "  -4.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.630"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9779"
$ws.Range("E35").Value = "  -2.92%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5897"
$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01594"
$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.983"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8577"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.026.70"
$ws.Range("E43").Value = "  -4.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.78"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.801.88"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("E46").Value = "  +3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.25"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.083"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.465"
$ws.Range("E51").Value = "  +0.37%  "
